# Generate Report for Handoff
# Rows 4-7 (314c8870..., 8f29e399..., b5477690..., ebb0f100...) move from
# "low" priority / "Ready for handoff" handoff status to "ht" priority with
# refreshed handoff timestamps, for both the zh-cn and de-de locale sheets
# (and the Overview roll-up sheet that mirrors the handoff-generated date).

$wb = $excel.ActiveWorkbook

# --- zh-cn sheet: rows 4-7 ---
$wsZh = $wb.Worksheets.Item("zh-cn")
foreach ($r in 4..7) {
    $wsZh.Cells.Item($r, 5).Value = "ht"                       # column E - Priority
    $wsZh.Cells.Item($r, 8).Value = "2016-09-04 04:34:58"       # column H - Latest Handoff Datetime
}

# --- de-de sheet: rows 4-7 ---
$wsDe = $wb.Worksheets.Item("de-de")
foreach ($r in 4..7) {
    $wsDe.Cells.Item($r, 5).Value = "ht"                        # column E - Priority
    $wsDe.Cells.Item($r, 8).Value = "2016-09-04 04:35:07"       # column H - Latest Handoff Datetime
}

# --- Overview sheet: rows 4-7, Latest HO Xliff Generate Date mirrors de-de's refreshed timestamp ---
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in 4..7) {
    $wsOverview.Cells.Item($r, 7).Value = "2016-09-04 04:35:07" # column G - Latest HO Xliff Generate Date
}
